$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header formatting: the "_old" / "_new" suffixes used for the two
#    compared-input-file column groups are replaced by the actual format
#    version names of those inputs ("_FV2410" / "_FV2504"), e.g.
#    "Segmentname_old" -> "Segmentname_FV2410" and
#    "Segmentname_new" -> "Segmentname_FV2504". The "diff" column in between
#    is left untouched.
# ---------------------------------------------------------------------------
$dim = $ws.UsedRange
$lastCol = $dim.Columns.Count

for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = $cell.Text

    if ($header -like "*_old") {
        $base = $header.Substring(0, $header.Length - [string]"_old".Length)
        $cell.Value = $base + "_FV2410"
    } elseif ($header -like "*_new") {
        $base = $header.Substring(0, $header.Length - [string]"_new".Length)
        $cell.Value = $base + "_FV2504"
    }
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row so it stays visible while scrolling.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3) Turn the whole used range into a proper Excel Table ("Table1") with
#    a header row and an autofilter, matching the sheet's data extent.
# ---------------------------------------------------------------------------
$lastRow = $dim.Rows.Count
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

Write-Output "Renamed headers, froze top row, and created $($tbl.Name) over $($tableRange.Address())"
